# Map108.xlsx v2 update: mirror column A into column B for every data row
# (duplicating each shared string), insert 6 new section-divider strings into
# the dialogue script, and retire the old D8 duplicate (EV025) in favor of B8.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the legacy duplicate in D8 (superseded by the new B-column mirroring)
$ws.Cells.Item(8,4).Value = $null

$t = '自動実行'
$ws.Cells.Item(1,1).Value = $t
$ws.Cells.Item(1,2).Value = $t

$t = 'いわ'
$ws.Cells.Item(2,1).Value = $t
$ws.Cells.Item(2,2).Value = $t

$t = ' <enemy:148>'
$ws.Cells.Item(3,1).Value = $t
$ws.Cells.Item(3,2).Value = $t

$t = '錆びついた鉄格子だ・・・'
$ws.Cells.Item(4,1).Value = $t
$ws.Cells.Item(4,2).Value = $t

$t = '食糧庫ドア'
$ws.Cells.Item(5,1).Value = $t
$ws.Cells.Item(5,2).Value = $t

$t = '鍵が刺さったままになっている・・・'
$ws.Cells.Item(6,1).Value = $t
$ws.Cells.Item(6,2).Value = $t

$t = '出現防止'
$ws.Cells.Item(7,1).Value = $t
$ws.Cells.Item(7,2).Value = $t

$t = 'EV025'
$ws.Cells.Item(8,1).Value = $t
$ws.Cells.Item(8,2).Value = $t

$t = 'リリーの日記'
$ws.Cells.Item(9,1).Value = $t
$ws.Cells.Item(9,2).Value = $t

$t = 'リリーの日記
ライム式ダイエット'
$ws.Cells.Item(10,1).Value = $t
$ws.Cells.Item(10,2).Value = $t

$t = '読んでみる'
$ws.Cells.Item(11,1).Value = $t
$ws.Cells.Item(11,2).Value = $t

$t = 'やめておく'
$ws.Cells.Item(12,1).Value = $t
$ws.Cells.Item(12,2).Value = $t

$t = 'EV054'
$ws.Cells.Item(13,1).Value = $t
$ws.Cells.Item(13,2).Value = $t

$t = '罠を踏んでしまった・・・！'
$ws.Cells.Item(14,1).Value = $t
$ws.Cells.Item(14,2).Value = $t

$t = 'なんと身体が小さくなってしまった・・・！'
$ws.Cells.Item(15,1).Value = $t
$ws.Cells.Item(15,2).Value = $t

$t = '--------キャラ指定-----座標設定------------'
$ws.Cells.Item(16,1).Value = $t

$t = '----------------------------------------'
$ws.Cells.Item(17,1).Value = $t

$t = '\n<リリー>なぁに？
また食べられたいの？ヘンタイ♥
わざと踏んだでしょそれ。'
$ws.Cells.Item(18,1).Value = $t
$ws.Cells.Item(18,2).Value = $t

$t = '\n<リリー>違うの？
本当に？ふーん。
ま、どっちでもいいけど。'
$ws.Cells.Item(19,1).Value = $t
$ws.Cells.Item(19,2).Value = $t

$t = 'ミニマムバトル'
$ws.Cells.Item(20,1).Value = $t
$ws.Cells.Item(20,2).Value = $t

$t = '\C[10]小説'
$ws.Cells.Item(21,1).Value = $t
$ws.Cells.Item(21,2).Value = $t

$t = '説明'
$ws.Cells.Item(22,1).Value = $t
$ws.Cells.Item(22,2).Value = $t

$t = '丸飲みサークルD-gateのCamelさんから頂いた
小説を元にしたイベントが見られます。
選択すると食べられてゲームオーバーになるので注意！'
$ws.Cells.Item(23,1).Value = $t
$ws.Cells.Item(23,2).Value = $t

$t = 'ポップアップ有効化'
$ws.Cells.Item(24,1).Value = $t
$ws.Cells.Item(24,2).Value = $t

$t = '\n<リリー>必死に抵抗しちゃって。
そんなに食べられたくなかった？
くすくす♥'
$ws.Cells.Item(25,1).Value = $t
$ws.Cells.Item(25,2).Value = $t

$t = '\n<リリー>そろそろ魔法が切れそうね。
じゃ、私のお口が恋しくなったらいつでも言うのよ？
また小さくして可愛がってあげる♥くすくす♥'
$ws.Cells.Item(26,1).Value = $t
$ws.Cells.Item(26,2).Value = $t

$t = '身体が元の大きさに戻った！'
$ws.Cells.Item(27,1).Value = $t
$ws.Cells.Item(27,2).Value = $t

$t = '報告書'
$ws.Cells.Item(28,1).Value = $t
$ws.Cells.Item(28,2).Value = $t

$t = 'お風呂場のモップ、ダメになってきちゃった！
根元がポキっと折れてただの『棒』になっちゃいそうなの。
アーコちゃん、もしモップあったらお願いね！'
$ws.Cells.Item(29,1).Value = $t
$ws.Cells.Item(29,2).Value = $t

$t = '\n<\n[1]>（淫魔の癖になんて庶民的な買い物・・・）'
$ws.Cells.Item(30,1).Value = $t
$ws.Cells.Item(30,2).Value = $t

$t = '\n<\n[1]>（モップ・・・武器になるかな？
こんな短い果物ナイフよりはマシだと思うけど・・・）'
$ws.Cells.Item(31,1).Value = $t
$ws.Cells.Item(31,2).Value = $t

$t = '\n<\n[1]>（モップか・・・武器に出来るか？
無いよりはましそうだけど・・・）'
$ws.Cells.Item(32,1).Value = $t
$ws.Cells.Item(32,2).Value = $t

$t = 'EV058'
$ws.Cells.Item(33,1).Value = $t
$ws.Cells.Item(33,2).Value = $t

$t = '何か使えそうなものは・・・'
$ws.Cells.Item(34,1).Value = $t
$ws.Cells.Item(34,2).Value = $t

$t = '豪華な装飾の付いた鍵だ・・・
どこの鍵だろう。'
$ws.Cells.Item(35,1).Value = $t
$ws.Cells.Item(35,2).Value = $t

$t = '\n<\n[3]>ごそごそして何か探してるのかな？
お腹空いたのー？'
$ws.Cells.Item(36,1).Value = $t
$ws.Cells.Item(36,2).Value = $t

$t = '\n<\n[3]>お腹空いたのなら私が何か作ってあげるよー♥
でもでもー、その前にー・・・'
$ws.Cells.Item(37,1).Value = $t
$ws.Cells.Item(37,2).Value = $t

$t = '\n<\n[3]>ぎゅーっ♥
私が先にご飯にするー♥'
$ws.Cells.Item(38,1).Value = $t
$ws.Cells.Item(38,2).Value = $t

$t = 'ーーーーー待機or攻撃アニメSE選択ーーーーー'
$ws.Cells.Item(39,1).Value = $t

$t = 'MP_SET_speed 5 50'
$ws.Cells.Item(40,1).Value = $t
$ws.Cells.Item(40,2).Value = $t

$t = '--------ここから誘惑------------'
$ws.Cells.Item(41,1).Value = $t

$t = '\n<\n[3]>わぁ♥
抱き着いたら乳首におちんちん入っちゃった♥
あははー♥'
$ws.Cells.Item(42,1).Value = $t
$ws.Cells.Item(42,2).Value = $t

$t = '\n<\n[3]>おっぱいで挟むのもいいけどー、
おっぱいの中に入っちゃうのもいいよねー♥
くすくす♥'
$ws.Cells.Item(43,1).Value = $t
$ws.Cells.Item(43,2).Value = $t

$t = '\n<\n[3]>ねーねー♥
腰動いちゃってるよー？
私の乳首でもぐもぐされるのどーお？気持ちいー？'
$ws.Cells.Item(44,1).Value = $t
$ws.Cells.Item(44,2).Value = $t

$t = '--------ここから選択肢-------------'
$ws.Cells.Item(45,1).Value = $t

$t = '\n<\n[3]>抜け出さなきゃなのに
にゅるにゅるおちんちんに絡みついてきて
逃げられないねー♥'
$ws.Cells.Item(46,1).Value = $t
$ws.Cells.Item(46,2).Value = $t

$t = 'MP_SET_speed 5 75'
$ws.Cells.Item(47,1).Value = $t
$ws.Cells.Item(47,2).Value = $t

$t = 'MP_SET_speed 5 100'
$ws.Cells.Item(48,1).Value = $t
$ws.Cells.Item(48,2).Value = $t

$t = '\n<\n[3]>もう離れられると思ったー？
ダメダメー♥
ほーらぴゅっぴゅしちゃいなよー♥'
$ws.Cells.Item(49,1).Value = $t
$ws.Cells.Item(49,2).Value = $t

$t = '\n<\n[3]>おいでおいでー♥'
$ws.Cells.Item(50,1).Value = $t
$ws.Cells.Item(50,2).Value = $t

$t = '--------ここからクリア------------'
$ws.Cells.Item(51,1).Value = $t

$t = '\n<\n[3]>むー・・・
カチカチにしてるのにまだ頑張るのー？そっかー。
悲しいなぁー・・・'
$ws.Cells.Item(52,1).Value = $t
$ws.Cells.Item(52,2).Value = $t

$t = '\n<\n[3]>でも、またしてほしくなったら言ってね？
いつでももぐもぐ食べてあげるから。
じゃ、また後でねー♥'
$ws.Cells.Item(53,1).Value = $t
$ws.Cells.Item(53,2).Value = $t

$t = 'アイテム'
$ws.Cells.Item(54,1).Value = $t
$ws.Cells.Item(54,2).Value = $t

$t = 'バルコニーの鍵を見つけた！
これで外に出られる・・・！'
$ws.Cells.Item(55,1).Value = $t
$ws.Cells.Item(55,2).Value = $t

$t = 'もう使えそうなものはない。'
$ws.Cells.Item(56,1).Value = $t
$ws.Cells.Item(56,2).Value = $t

$t = '食材が並んでいる・・・'
$ws.Cells.Item(57,1).Value = $t
$ws.Cells.Item(57,2).Value = $t

$t = '\n<\n[1]>（淫魔も普通のご飯食べるのかな。
じゃあそれだけ食べてればいいのに・・・）'
$ws.Cells.Item(58,1).Value = $t
$ws.Cells.Item(58,2).Value = $t

$t = '\n<\n[1]>（食べるものには困らないな・・・）'
$ws.Cells.Item(59,1).Value = $t
$ws.Cells.Item(59,2).Value = $t

$t = '時計が時を刻んでいる・・・'
$ws.Cells.Item(60,1).Value = $t
$ws.Cells.Item(60,2).Value = $t

$t = '\n<\n[1]>（どれだけの時間、ここにいるんだろう・・・）'
$ws.Cells.Item(61,1).Value = $t
$ws.Cells.Item(61,2).Value = $t

$t = '\n<\n[1]>（時間の感覚がおかしくなる・・・
早くここから抜け出さないと・・・）'
$ws.Cells.Item(62,1).Value = $t
$ws.Cells.Item(62,2).Value = $t

$t = 'リリーの日記
手配書を出そう'
$ws.Cells.Item(63,1).Value = $t
$ws.Cells.Item(63,2).Value = $t

$t = '食料保管庫と書かれている・・・'
$ws.Cells.Item(64,1).Value = $t
$ws.Cells.Item(64,2).Value = $t

$t = '\n<\n[1]>（何か甘い物、ないかな？）'
$ws.Cells.Item(65,1).Value = $t
$ws.Cells.Item(65,2).Value = $t

$t = '\n<\n[1]>（腐ってなければいいが・・・）'
$ws.Cells.Item(66,1).Value = $t
$ws.Cells.Item(66,2).Value = $t

$t = '来客に食事を提供する際の注意
1：身だしなみを整える
2：おもてなしの心をもつ'
$ws.Cells.Item(67,1).Value = $t
$ws.Cells.Item(67,2).Value = $t

$t = '3：適温で提供する
4：なるべくお酒を勧める（アルコール度数が高い物）
5：無駄な話をしない'
$ws.Cells.Item(68,1).Value = $t
$ws.Cells.Item(68,2).Value = $t

$t = '6：態度が怪しい場合はガード、またはメイド長へ連絡
7：問題ない場合は笑顔で帰ってもらう'
$ws.Cells.Item(69,1).Value = $t
$ws.Cells.Item(69,2).Value = $t

$t = 'EV074'
$ws.Cells.Item(70,1).Value = $t
$ws.Cells.Item(70,2).Value = $t

$t = 'EV075'
$ws.Cells.Item(71,1).Value = $t
$ws.Cells.Item(71,2).Value = $t

$t = '錆びた鉄格子を開けるには・・・'
$ws.Cells.Item(72,1).Value = $t
$ws.Cells.Item(72,2).Value = $t

$t = '答えを見る'
$ws.Cells.Item(73,1).Value = $t
$ws.Cells.Item(73,2).Value = $t

$t = '見ない'
$ws.Cells.Item(74,1).Value = $t
$ws.Cells.Item(74,2).Value = $t

$t = '武器で攻撃してください。'
$ws.Cells.Item(75,1).Value = $t
$ws.Cells.Item(75,2).Value = $t

$t = 'ビンにはピクルスが入っている。'
$ws.Cells.Item(76,1).Value = $t
$ws.Cells.Item(76,2).Value = $t

$t = '\n<\n[1]>（ハンバーガーに挟むと美味しい。）'
$ws.Cells.Item(77,1).Value = $t
$ws.Cells.Item(77,2).Value = $t

$t = '\n<\n[1]>（こんな酸っぱいものは
ハンバーガーに挟むべきではない。）'
$ws.Cells.Item(78,1).Value = $t
$ws.Cells.Item(78,2).Value = $t

$t = '様々な種類のビネガーが置かれている。'
$ws.Cells.Item(79,1).Value = $t
$ws.Cells.Item(79,2).Value = $t

$t = '\n<\n[1]>（お酢を飲むと体が柔らかくなるってお母さんが言ってた。
多分気のせい。）'
$ws.Cells.Item(80,1).Value = $t
$ws.Cells.Item(80,2).Value = $t

$t = '\n<\n[1]>（これが扱える人は自炊のレベルが高いと思う。）'
$ws.Cells.Item(81,1).Value = $t
$ws.Cells.Item(81,2).Value = $t

$t = '料理酒が並んでいる。'
$ws.Cells.Item(82,1).Value = $t
$ws.Cells.Item(82,2).Value = $t

$t = '\n<\n[1]>（お肉の臭みを取ったり柔らかくしたり出来る。
今要るかと言われると要らない。）'
$ws.Cells.Item(83,1).Value = $t
$ws.Cells.Item(83,2).Value = $t

$t = '\n<\n[1]>（何のために料理に入れるのかさっぱり分からない。）'
$ws.Cells.Item(84,1).Value = $t
$ws.Cells.Item(84,2).Value = $t

$t = 'まな板。'
$ws.Cells.Item(85,1).Value = $t
$ws.Cells.Item(85,2).Value = $t

$t = '\n<\n[1]>（淫魔の出してくる食べ物なんてと思っていたけど
無駄に美味しいんだよなぁ・・・
逆になんかやだ。）'
$ws.Cells.Item(86,1).Value = $t
$ws.Cells.Item(86,2).Value = $t

$t = '\n<\n[1]>（ここの食事が妙に美味しいのも
何か逆に良くない気がする・・・）'
$ws.Cells.Item(87,1).Value = $t
$ws.Cells.Item(87,2).Value = $t

$t = 'EV080'
$ws.Cells.Item(88,1).Value = $t
$ws.Cells.Item(88,2).Value = $t

$t = 'バルコニーの鍵！
これで外に出られそうだ・・・！'
$ws.Cells.Item(89,1).Value = $t
$ws.Cells.Item(89,2).Value = $t

$t = '身体が小さくなってしまった・・・！'
$ws.Cells.Item(90,1).Value = $t
$ws.Cells.Item(90,2).Value = $t

$t = '\n<リリー>あはは！ひっかかった！
可愛らしいサイズになったわねぇ♥くすくす♥'
$ws.Cells.Item(91,1).Value = $t
$ws.Cells.Item(91,2).Value = $t

$t = '\n<リリー>あんた食糧庫で何してたの？
お腹空いた？'
$ws.Cells.Item(92,1).Value = $t
$ws.Cells.Item(92,2).Value = $t

$t = '\n<リリー>でもざーんねん。
食事するのは私よ。'
$ws.Cells.Item(93,1).Value = $t
$ws.Cells.Item(93,2).Value = $t

$t = '\n<リリー>そろそろ魔法が切れそうね。
じゃ、私のお口が恋しくなったらいつでも言うのよ？
また可愛がってあげる♥くすくす♥'
$ws.Cells.Item(94,1).Value = $t
$ws.Cells.Item(94,2).Value = $t

$t = '炎の魔導書・・・
装備すれば何かを燃やすことが出来そうだ。'
$ws.Cells.Item(95,1).Value = $t
$ws.Cells.Item(95,2).Value = $t

$t = '\n<リリー>捕まえたっ♥
ほらこっち向きなさい♥'
$ws.Cells.Item(96,1).Value = $t
$ws.Cells.Item(96,2).Value = $t

$t = '\n<\n[3]>ちゅー♥
あっは♥いきなりキスしちゃった♥
くすくす♥'
$ws.Cells.Item(97,1).Value = $t
$ws.Cells.Item(97,2).Value = $t

$t = '\n<\n[3]>あんたってマゾだからさー、
こうやってキスされるとすぐスイッチ入っちゃうでしょ。
ちゅっ♥ちゅっ♥'
$ws.Cells.Item(98,1).Value = $t
$ws.Cells.Item(98,2).Value = $t

$t = '\n<\n[3]>ほらビンビンに勃起した。あーぁ。
ここに挟まれたらもう逃げられないわね。'
$ws.Cells.Item(99,1).Value = $t
$ws.Cells.Item(99,2).Value = $t

$t = '\n<\n[3]>あら、振り払おうとしてるの？
全然力が入ってないけど♥
おちんちんは逃げたくないみたいよ？'
$ws.Cells.Item(100,1).Value = $t
$ws.Cells.Item(100,2).Value = $t

$t = '\n<\n[3]>イキなさい。
あんたの負けよ。
くすくす♥'
$ws.Cells.Item(101,1).Value = $t
$ws.Cells.Item(101,2).Value = $t

$t = '\n<リリー>なにあんた。必死になっちゃって。
分かったわよ。
離れてあげる。'
$ws.Cells.Item(102,1).Value = $t
$ws.Cells.Item(102,2).Value = $t

$t = '\n<リリー>あらあら大変。
我慢汁だらだらじゃないの。
床汚さないでよ。'
$ws.Cells.Item(103,1).Value = $t
$ws.Cells.Item(103,2).Value = $t

$t = 'EV081'
$ws.Cells.Item(104,1).Value = $t
$ws.Cells.Item(104,2).Value = $t

$t = 'EV082'
$ws.Cells.Item(105,1).Value = $t
$ws.Cells.Item(105,2).Value = $t
